$wb = $excel.ActiveWorkbook

# Both the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet
# contain the same data table and need the same updates to column F
# ("想去人数" / number interested), reflecting refreshed counts from a
# newer generation of the page (gh-pages output update).
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 2433
    $ws.Range("F4").Value = 455
    $ws.Range("F6").Value = 6517
    $ws.Range("F7").Value = 353
    $ws.Range("F8").Value = 131
}
